$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Freq. of Occurrence")
$ws.Activate()

# Update the color value for row 121 (climatecolor column F)
$ws.Range("F121").Value = "#F0E442"

# Move the view: scroll so A105 is the top-left visible cell, and select G120
$excel.ActiveWindow.ScrollRow = 105
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("G120").Select()
